$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.105.45"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.867.86"
$ws.Range("E3").Value = "  +3.74%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.66"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4993"
$ws.Range("E7").Value = "  -2.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3919"
$ws.Range("E8").Value = "  +1.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09733"
$ws.Range("E9").Value = "  +25.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.131"
$ws.Range("E10").Value = "  +2.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "40.96"
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.463"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.93"
$ws.Range("E13").Value = "  +3.30%  "
$ws.Range("D14").Value = "1.871.77"
$ws.Range("E14").Value = "  +4.18%  "
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.382"
$ws.Range("E16").Value = "  +1.45%  "
$ws.Range("E17").Value = "  +4.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.05"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.44"
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.120"
$ws.Range("E22").Value = "  +2.95%  "
$ws.Range("D23").Value = "28.173.66"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.32"
$ws.Range("E24").Value = "  +2.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.277"
$ws.Range("E25").Value = "  +1.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.559"
$ws.Range("E26").Value = "  +5.72%  "
$ws.Range("D27").Value = "2.080.88"
$ws.Range("E27").Value = "  +3.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.07"
$ws.Range("E28").Value = "  +4.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "157.63"
$ws.Range("E29").Value = "  -1.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.32"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("E31").Value = "  -2.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.060"
$ws.Range("E32").Value = "  +1.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.610"
$ws.Range("E33").Value = "  +1.43%  "
$ws.Range("E34").Value = "  -0.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06754"
$ws.Range("E35").Value = "  -3.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.454"
$ws.Range("E36").Value = "  +4.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02384"
$ws.Range("E37").Value = "  +1.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2181"
$ws.Range("E38").Value = "  +0.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.47"
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.011"
$ws.Range("E40").Value = "  -0.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6283"
$ws.Range("E41").Value = "  +2.77%  "
$ws.Range("E42").Value = "  +1.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.54"
$ws.Range("E44").Value = "  +3.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5986"
$ws.Range("E45").Value = "  +1.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.658"
$ws.Range("E46").Value = "  -1.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.267"
$ws.Range("E47").Value = "  -2.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.23"
$ws.Range("E48").Value = "  -0.78%  "
$ws.Range("E49").Value = "  +3.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.195"
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06835"
$ws.Range("E51").Value = "  +1.53%  "
